$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2989916666666667
$ws.Range("H2").Value = 0.896975
$ws.Range("I2").Value = 0.01120651476222736
$ws.Range("J2").Value = 0.01120651476222736
$ws.Range("M2").Value = 32.51511900000001
$ws.Range("N2").Value = 97.54535700000001
$ws.Range("O2").Value = 0.218203973858649
$ws.Range("P2").Value = 0.2182039738586489
$ws.Range("Q2").Value = 9.721749621675002
$ws.Range("R2").Value = 87.495746595075
$ws.Range("S2").Value = 0.002445306054223623
$ws.Range("T2").Value = 0.002445306054223622
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2989916666666667
$ws.Range("H3").Value = 0.896975
$ws.Range("I3").Value = 0.01120651476222736
$ws.Range("J3").Value = 0.01120651476222736
$ws.Range("O3").Value = 0.6017421411306194
$ws.Range("P3").Value = 0.6017421411306194
$ws.Range("Q3").Value = 26.809715375175
$ws.Range("R3").Value = 241.287438376575
$ws.Range("S3").Value = 0.006743432187634586
$ws.Range("T3").Value = 0.006743432187634586
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2989916666666667
$ws.Range("H4").Value = 0.896975
$ws.Range("I4").Value = 0.01120651476222736
$ws.Range("J4").Value = 0.01120651476222736
$ws.Range("M4").Value = 26.830279
$ws.Range("N4").Value = 80.490837
$ws.Range("O4").Value = 0.1800538850107317
$ws.Range("P4").Value = 0.1800538850107317
$ws.Range("Q4").Value = 8.022029835341666
$ws.Range("R4").Value = 72.198268518075
$ws.Range("S4").Value = 0.002017776520369152
$ws.Range("T4").Value = 0.002017776520369152
# Row 5
$ws.Range("I5").Value = 0.9329357354307251
$ws.Range("J5").Value = 0.9329357354307249
$ws.Range("M5").Value = 32.51511900000001
$ws.Range("N5").Value = 97.54535700000001
$ws.Range("O5").Value = 0.218203973858649
$ws.Range("P5").Value = 0.2182039738586489
$ws.Range("Q5").Value = 809.3299143763471
$ws.Range("R5").Value = 7283.969229387124
$ws.Range("S5").Value = 0.2035702848257254
$ws.Range("T5").Value = 0.2035702848257253
# Row 6
$ws.Range("I6").Value = 0.9329357354307251
$ws.Range("J6").Value = 0.9329357354307249
$ws.Range("O6").Value = 0.6017421411306194
$ws.Range("P6").Value = 0.6017421411306194
$ws.Range("S6").Value = 0.5613867469753536
$ws.Range("T6").Value = 0.5613867469753535
# Row 7
$ws.Range("I7").Value = 0.9329357354307251
$ws.Range("J7").Value = 0.9329357354307249
$ws.Range("M7").Value = 26.830279
$ws.Range("N7").Value = 80.490837
$ws.Range("O7").Value = 0.1800538850107317
$ws.Range("P7").Value = 0.1800538850107317
$ws.Range("Q7").Value = 667.8292460120937
$ws.Range("R7").Value = 6010.463214108843
$ws.Range("S7").Value = 0.1679787036296462
$ws.Range("T7").Value = 0.1679787036296461
# Row 8
$ws.Range("G8").Value = 1.490294
$ws.Range("H8").Value = 4.470882
$ws.Range("I8").Value = 0.05585774980704767
$ws.Range("J8").Value = 0.05585774980704766
$ws.Range("M8").Value = 32.51511900000001
$ws.Range("N8").Value = 97.54535700000001
$ws.Range("O8").Value = 0.218203973858649
$ws.Range("P8").Value = 0.2182039738586489
$ws.Range("Q8").Value = 48.457086754986
$ws.Range("R8").Value = 436.113780794874
$ws.Range("S8").Value = 0.01218838297869998
$ws.Range("T8").Value = 0.01218838297869998
# Row 9
$ws.Range("G9").Value = 1.490294
$ws.Range("H9").Value = 4.470882
$ws.Range("I9").Value = 0.05585774980704767
$ws.Range("J9").Value = 0.05585774980704766
$ws.Range("O9").Value = 0.6017421411306194
$ws.Range("P9").Value = 0.6017421411306194
$ws.Range("Q9").Value = 133.630339637106
$ws.Range("R9").Value = 1202.673056733954
$ws.Range("S9").Value = 0.03361196196763131
$ws.Range("T9").Value = 0.0336119619676313
# Row 10
$ws.Range("G10").Value = 1.490294
$ws.Range("H10").Value = 4.470882
$ws.Range("I10").Value = 0.05585774980704767
$ws.Range("J10").Value = 0.05585774980704766
$ws.Range("M10").Value = 26.830279
$ws.Range("N10").Value = 80.490837
$ws.Range("O10").Value = 0.1800538850107317
$ws.Range("P10").Value = 0.1800538850107317
$ws.Range("Q10").Value = 39.98500381202599
$ws.Range("R10").Value = 359.865034308234
$ws.Range("S10").Value = 0.01005740486071638
$ws.Range("T10").Value = 0.01005740486071638
